$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.621.47'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.852.41'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '243.67'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '0.6520'
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '48.25'
$ws.Range("E8").Value = '  +4.18%  '
$ws.Range("D9").Value = '0.07495'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '0.2975'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '24.52'
$ws.Range("E11").Value = '  +2.44%  '
$ws.Range("D12").Value = '0.07641'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = '1.861.75'
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '5.044'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").Value = '0.6863'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").Value = '83.57'
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("D17").Value = '0.000009618'
$ws.Range("E17").Value = '  +2.62%  '
$ws.Range("D18").Value = '6.147'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").Value = '29.664.99'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '2.112.14'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("D21").Value = '237.14'
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").Value = '12.63'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").Value = '7.710'
$ws.Range("E24").Value = '  +4.30%  '
$ws.Range("D25").Value = '1.002'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '158.20'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").Value = '0.1425'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").Value = '8.527'
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("D29").Value = '17.86'
$ws.Range("E29").Value = '  -0.45%  '
$ws.Range("D30").Value = '0.06043'
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '1.488'
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").Value = '1.265'
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = '4.146'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '4.078'
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("D35").Value = '1.187'
$ws.Range("E35").Value = '  +3.01%  '
$ws.Range("D36").Value = '1.872'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").Value = '0.7282'
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").Value = '2.600'
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("D39").Value = '2.800'
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("D40").Value = '0.01786'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").Value = '1.203.12'
$ws.Range("E41").Value = '  -2.44%  '
$ws.Range("D42").Value = '6.297'
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '0.9111'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '2.024.63'
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("D46").Value = '101.31'
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").Value = '66.63'
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '0.00000000124'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").Value = '7.351'
$ws.Range("E49").Value = '  +9.54%  '
$ws.Range("D50").Value = '0.4058'
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("D51").Value = '9.167'
$ws.Range("E51").Value = '  -1.21%  '
